$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header style from E1 to F1 so it reuses the same cellXf (bold/border/alignment)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Populate time_taken values for each data row
$ws.Range("F2").Value = "2021-10-05 13:39:33.453154"
$ws.Range("F3").Value = "2021-10-05 13:39:33.453167"
$ws.Range("F4").Value = "2021-10-05 13:39:33.453171"
$ws.Range("F5").Value = "2021-10-05 13:39:33.453174"
$ws.Range("F6").Value = "2021-10-05 13:39:33.453178"
$ws.Range("F7").Value = "2021-10-05 13:39:33.453181"
$ws.Range("F8").Value = "2021-10-05 13:39:33.453184"
$ws.Range("F9").Value = "2021-10-05 13:39:33.453187"
$ws.Range("F10").Value = "2021-10-05 13:39:33.453190"
$ws.Range("F11").Value = "2021-10-05 13:39:33.453193"
$ws.Range("F12").Value = "2021-10-05 13:39:33.453196"
$ws.Range("F13").Value = "2021-10-05 13:39:33.453199"
$ws.Range("F14").Value = "2021-10-05 13:39:33.453202"
$ws.Range("F15").Value = "2021-10-05 13:39:33.453205"
$ws.Range("F16").Value = "2021-10-05 13:39:33.453208"
$ws.Range("F17").Value = "2021-10-05 13:39:33.453211"
$ws.Range("F18").Value = "2021-10-05 13:39:33.453215"
$ws.Range("F19").Value = "2021-10-05 13:39:33.453218"
$ws.Range("F20").Value = "2021-10-05 13:39:33.453221"
$ws.Range("F21").Value = "2021-10-05 13:39:33.453224"
$ws.Range("F22").Value = "2021-10-05 13:39:33.453227"
$ws.Range("F23").Value = "2021-10-05 13:39:33.453229"
$ws.Range("F24").Value = "2021-10-05 13:39:33.453232"
$ws.Range("F25").Value = "2021-10-05 13:39:33.453235"
$ws.Range("F26").Value = "2021-10-05 13:39:33.453239"
$ws.Range("F27").Value = "2021-10-05 13:39:33.453242"
$ws.Range("F28").Value = "2021-10-05 13:39:33.453245"
$ws.Range("F29").Value = "2021-10-05 13:39:33.453248"
$ws.Range("F30").Value = "2021-10-05 13:39:33.453251"
$ws.Range("F31").Value = "2021-10-05 13:39:33.453254"
$ws.Range("F32").Value = "2021-10-05 13:39:33.453257"
$ws.Range("F33").Value = "2021-10-05 13:39:33.453260"
$ws.Range("F34").Value = "2021-10-05 13:39:33.453263"
$ws.Range("F35").Value = "2021-10-05 13:39:33.453266"
$ws.Range("F36").Value = "2021-10-05 13:39:33.453269"
$ws.Range("F37").Value = "2021-10-05 13:39:33.453272"
$ws.Range("F38").Value = "2021-10-05 13:39:33.453275"
$ws.Range("F39").Value = "2021-10-05 13:39:33.453278"
$ws.Range("F40").Value = "2021-10-05 13:39:33.453281"
$ws.Range("F41").Value = "2021-10-05 13:39:33.453284"
$ws.Range("F42").Value = "2021-10-05 13:39:33.453288"
$ws.Range("F43").Value = "2021-10-05 13:39:33.453291"
$ws.Range("F44").Value = "2021-10-05 13:39:33.453294"
$ws.Range("F45").Value = "2021-10-05 13:39:33.453297"
$ws.Range("F46").Value = "2021-10-05 13:39:33.453300"
$ws.Range("F47").Value = "2021-10-05 13:39:33.453303"
$ws.Range("F48").Value = "2021-10-05 13:39:33.453306"
$ws.Range("F49").Value = "2021-10-05 13:39:33.453309"
$ws.Range("F50").Value = "2021-10-05 13:39:33.453312"
$ws.Range("F51").Value = "2021-10-05 13:39:33.453315"
$ws.Range("F52").Value = "2021-10-05 13:39:33.453318"
$ws.Range("F53").Value = "2021-10-05 13:39:33.453321"
$ws.Range("F54").Value = "2021-10-05 13:39:33.453325"
$ws.Range("F55").Value = "2021-10-05 13:39:33.453328"
$ws.Range("F56").Value = "2021-10-05 13:39:33.453331"
$ws.Range("F57").Value = "2021-10-05 13:39:33.453334"
$ws.Range("F58").Value = "2021-10-05 13:39:33.453337"
$ws.Range("F59").Value = "2021-10-05 13:39:33.453339"
$ws.Range("F60").Value = "2021-10-05 13:39:33.453342"
$ws.Range("F61").Value = "2021-10-05 13:39:33.453345"
$ws.Range("F62").Value = "2021-10-05 13:39:33.453348"
$ws.Range("F63").Value = "2021-10-05 13:39:33.453351"
$ws.Range("F64").Value = "2021-10-05 13:39:33.453354"
$ws.Range("F65").Value = "2021-10-05 13:39:33.453357"
$ws.Range("F66").Value = "2021-10-05 13:39:33.453361"
$ws.Range("F67").Value = "2021-10-05 13:39:33.453365"
$ws.Range("F68").Value = "2021-10-05 13:39:33.453368"
$ws.Range("F69").Value = "2021-10-05 13:39:33.453371"
$ws.Range("F70").Value = "2021-10-05 13:39:33.453374"
$ws.Range("F71").Value = "2021-10-05 13:39:33.453378"
$ws.Range("F72").Value = "2021-10-05 13:39:33.453381"
$ws.Range("F73").Value = "2021-10-05 13:39:33.453384"
$ws.Range("F74").Value = "2021-10-05 13:39:33.453387"
$ws.Range("F75").Value = "2021-10-05 13:39:33.453390"
$ws.Range("F76").Value = "2021-10-05 13:39:33.453393"
$ws.Range("F77").Value = "2021-10-05 13:39:33.453396"
